$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws 'D2' '28.243.84'
Set-TextValue $ws 'E2' '  -0.57%  '
Set-TextValue $ws 'D3' '1.805.51'
Set-TextValue $ws 'D4' '1.003'
Set-TextValue $ws 'E4' '  +0.15%  '
Set-TextValue $ws 'D5' '314.74'
Set-TextValue $ws 'E5' '  -0.27%  '
Set-TextValue $ws 'D6' '1.002'
Set-TextValue $ws 'E6' '  +0.05%  '
Set-TextValue $ws 'D7' '0.5253'
Set-TextValue $ws 'E7' '  +2.30%  '
Set-TextValue $ws 'D8' '0.3834'
Set-TextValue $ws 'E8' '  -2.15%  '
Set-TextValue $ws 'D9' '0.08010'
Set-TextValue $ws 'E9' '  +1.85%  '
Set-TextValue $ws 'B10' 'OKB'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws 'D10' '41.47'
Set-TextValue $ws 'E10' '  -0.71%  '
Set-TextValue $ws 'B11' 'Polygon'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws 'D11' '1.104'
Set-TextValue $ws 'E11' '  -0.43%  '
Set-TextValue $ws 'D12' '6.358'
Set-TextValue $ws 'E12' '  +1.87%  '
Set-TextValue $ws 'E13' '  +0.11%  '
Set-TextValue $ws 'D14' '20.65'
Set-TextValue $ws 'E14' '  -1.43%  '
Set-TextValue $ws 'E15' '  -1.45%  '
Set-TextValue $ws 'D16' '1.802.29'
Set-TextValue $ws 'E16' '  -1.08%  '
Set-TextValue $ws 'D17' '92.74'
Set-TextValue $ws 'E17' '  +0.24%  '
Set-TextValue $ws 'D18' '0.00001100'
Set-TextValue $ws 'E18' '  -2.39%  '
Set-TextValue $ws 'D19' '0.06609'
Set-TextValue $ws 'E19' '  -0.21%  '
Set-TextValue $ws 'D20' '1.001'
Set-TextValue $ws 'E20' '  +0.03%  '
Set-TextValue $ws 'D21' '17.41'
Set-TextValue $ws 'E21' '  -1.49%  '
Set-TextValue $ws 'D22' '5.982'
Set-TextValue $ws 'E22' '  -1.63%  '
Set-TextValue $ws 'D23' '28.311.76'
Set-TextValue $ws 'E23' '  -0.42%  '
Set-TextValue $ws 'D24' '11.21'
Set-TextValue $ws 'E24' '  -0.41%  '
Set-TextValue $ws 'D25' '2.235'
Set-TextValue $ws 'E25' '  -1.37%  '
Set-TextValue $ws 'D26' '159.92'
Set-TextValue $ws 'E26' '  +3.64%  '
Set-TextValue $ws 'D27' '20.53'
Set-TextValue $ws 'E27' '  -2.61%  '
Set-TextValue $ws 'D28' '2.010.56'
Set-TextValue $ws 'E28' '  -0.86%  '
Set-TextValue $ws 'D29' '2.379'
Set-TextValue $ws 'E29' '  -0.70%  '
Set-TextValue $ws 'D30' '123.15'
Set-TextValue $ws 'E30' '  -1.77%  '
Set-TextValue $ws 'D31' '0.1084'
Set-TextValue $ws 'E31' '  -1.46%  '
Set-TextValue $ws 'D32' '1.060'
Set-TextValue $ws 'E32' '  -3.77%  '
Set-TextValue $ws 'D33' '3.674'
Set-TextValue $ws 'E33' '  +0.73%  '
Set-TextValue $ws 'D34' '5.573'
Set-TextValue $ws 'E34' '  -1.70%  '
Set-TextValue $ws 'D35' '0.07271'
Set-TextValue $ws 'E35' '  +3.31%  '
Set-TextValue $ws 'D36' '12.46'
Set-TextValue $ws 'B37' 'FraxShare'
Set-TextValue $ws 'C37' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D37' '8.917'
Set-TextValue $ws 'E37' '  +1.51%  '
Set-TextValue $ws 'B38' 'Algorand'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D38' '0.2173'
Set-TextValue $ws 'E38' '  -1.68%  '
Set-TextValue $ws 'B39' 'VeChain'
Set-TextValue $ws 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D39' '0.02322'
Set-TextValue $ws 'E39' '  +0.05%  '
Set-TextValue $ws 'D40' '5.150'
Set-TextValue $ws 'E40' '  -0.80%  '
Set-TextValue $ws 'D41' '0.6237'
Set-TextValue $ws 'E41' '  -0.22%  '
Set-TextValue $ws 'E42' '  -0.73%  '
Set-TextValue $ws 'D43' '1.373'
Set-TextValue $ws 'E43' '  -1.30%  '
Set-TextValue $ws 'B44' 'Decentraland'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws 'D44' '0.6039'
Set-TextValue $ws 'E44' '  +2.72%  '
Set-TextValue $ws 'B45' 'EnergySwap'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D45' '13.22'
Set-TextValue $ws 'E45' '  -1.31%  '
Set-TextValue $ws 'D46' '3.774'
Set-TextValue $ws 'E46' '  +1.07%  '
Set-TextValue $ws 'D47' '127.14'
Set-TextValue $ws 'E47' '  +2.06%  '
Set-TextValue $ws 'D48' '1.213'
Set-TextValue $ws 'E48' '  +1.74%  '
Set-TextValue $ws 'D49' '1.938'
Set-TextValue $ws 'E49' '  -1.58%  '
Set-TextValue $ws 'D50' '0.06846'
Set-TextValue $ws 'E50' '  -0.66%  '
Set-TextValue $ws 'D51' '73.29'
Set-TextValue $ws 'E51' '  -1.16%  '
